$wb = $excel.ActiveWorkbook

# Update "Full results" sheet (sheet1)
$ws1 = $wb.Worksheets.Item("Full results")

# Row 2: income / NULL MODEL
$ws1.Range("C2").Value = 0.97278525902507
$ws1.Range("D2").Value = 0.0272435083185367
$ws1.Range("E2").Value = 1.00002876734361
$ws1.Range("J2").Value = 0.0272427246177168
$ws1.Range("K2").Value = 0.0202315568872732
$ws1.Range("L2").Value = 0.185714405750598
$ws1.Range("M2").Value = 0.177081635541021
$ws1.Range("N2").Value = 0.205945962637871

# Row 3: income / CONDITIONAL MODEL
$ws1.Range("F3").Value = 0.981418277576514
$ws1.Range("G3").Value = 0.0202321388954219

# Row 4: income / COMPLETE MODEL
$ws1.Range("H4").Value = 0.795698529315793
$ws1.Range("I4").Value = 0.00682349204130963
$ws1.Range("O4").Value = 0.204324360158738

# Update "For plotting" sheet (sheet2)
$ws2 = $wb.Worksheets.Item("For plotting")

# Row 2: Sibcorr / income
$ws2.Range("C2").Value = 0.0272427246177168
$ws2.Range("D2").Value = -0.0208666986712964
$ws2.Range("E2").Value = 0.07535214790673

# Row 3: IOLIB / income
$ws2.Range("C3").Value = 0.205945962637871
$ws2.Range("D3").Value = 0.157339772481215
$ws2.Range("E3").Value = 0.254552152794528

# Row 4: IORAD / income
$ws2.Range("C4").Value = 0.204324360158738
$ws2.Range("D4").Value = 0.158648827976346
$ws2.Range("E4").Value = 0.24999989234113
